$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scraper was re-run against the source and now also reports each
# player's "height" and "weight". Those two new stat columns land right
# before the existing "fantasy points" column, so inserting them at E:F
# pushes the old "fantasy points" column from E out to G.
$ws.Range("E1:F1").EntireColumn.Insert()

# New header cells for the inserted columns.
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"

# Give the two new header cells the same formatting (bold, bordered,
# centered) as the rest of the header row.
$ws.Range("B1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Every player row in this sheet is the same person (Vance McDonald), so
# height/weight are constant across all the data rows.
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.333333333333333
    $ws.Cells.Item($r, 6).Value = 267
}
